$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.457.82'
$ws.Range("E2").Value = '  -2.70%  '
$ws.Range("D3").Value = '2.229.55'
$ws.Range("E3").Value = '  -2.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '111.18'
$ws.Range("E5").Value = '  -8.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '296.78'
$ws.Range("E6").Value = '  +11.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.625'
$ws.Range("E7").Value = '  -3.81%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.610'
$ws.Range("E9").Value = '  -3.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '44.49'
$ws.Range("E10").Value = '  -8.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0918'
$ws.Range("E11").Value = '  -3.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.26'
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.83'
$ws.Range("E13").Value = '  -5.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.01'
$ws.Range("E14").Value = '  +10.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.103'
$ws.Range("E15").Value = '  -2.79%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.10'
$ws.Range("E16").Value = '  -3.80%  '
$ws.Range("D17").Value = '2.561.12'
$ws.Range("D18").Value = '2.226.25'
$ws.Range("E18").Value = '  -2.23%  '
$ws.Range("D19").Value = '42.497.35'
$ws.Range("E19").Value = '  -2.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.35'
$ws.Range("E20").Value = '  +5.63%  '
$ws.Range("E21").Value = '  -4.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.88'
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.52'
$ws.Range("E23").Value = '  +21.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.39'
$ws.Range("E24").Value = '  -1.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '229.99'
$ws.Range("E25").Value = '  -2.54%  '
$ws.Range("E26").Value = '  -3.54%  '
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -1.66%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.62'
$ws.Range("E28").Value = '  -3.16%  '
$ws.Range("E29").Value = '  -1.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.54'
$ws.Range("E30").Value = '  -11.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.25'
$ws.Range("E31").Value = '  -4.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '173.36'
$ws.Range("E32").Value = '  -0.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.05'
$ws.Range("E33").Value = '  -3.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0898'
$ws.Range("E34").Value = '  -3.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.78'
$ws.Range("E35").Value = '  -0.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.13'
$ws.Range("E36").Value = '  +10.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.33'
$ws.Range("E37").Value = '  +0.64%  '
$ws.Range("E38").Value = '  -3.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0378'
$ws.Range("E39").Value = '  -2.03%  '
$ws.Range("E40").Value = '  -5.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.41'
$ws.Range("E41").Value = '  -5.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.47'
$ws.Range("E42").Value = '  -2.24%  '
$ws.Range("E43").Value = '  -2.09%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.35%  '
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.75'
$ws.Range("E45").Value = '  -7.00%  '
$ws.Range("E46").Value = '  -4.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.44'
$ws.Range("E47").Value = '  -8.39%  '
$ws.Range("E48").Value = '  +3.59%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.59'
$ws.Range("E49").Value = '  +0.92%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.18'
$ws.Range("E50").Value = '  -0.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.64'
$ws.Range("E51").Value = '  +6.35%  '
